# Unify the contact name in both rows and refresh the e-mail addresses,
# turning the e-mail cells into real hyperlinks (mailto:) as in the
# updated "contacts" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Names: use the same ("unique") name for both rows ---
$ws.Range("A1").Value = "FirstUser"
$ws.Range("A2").Value = "FirstUser"

# --- E-mails: replace the old addresses with the new ones ---
$ws.Range("B1").Value = "aa@mail.com"
$ws.Range("B2").Value = "asd@Gmail.com"

# --- Turn the e-mail cells into clickable mailto hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:aa@mail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:asd@Gmail.com")

# Re-apply the workbook's built-in hyperlink cell style so both cells
# keep using the same style already present in the file (instead of a
# freshly minted one).
$ws.Range("B1").Style = "Гиперссылка"
$ws.Range("B2").Style = "Гиперссылка"
